$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-300).
# All of these were bumped from 45182 (2023-09-13) to 45184 (2023-09-15).
$ws.Range("C2:C300").Value = 45184
